# Updated symbol list on Fri Jan 20 04:57:43 UTC 2023 with GitHub Actions
# Refresh Price (D) and Volume(1h) (E) columns for the crypto ranking rows
# that moved since the previous scrape. Values are plain text (not real
# numbers) in this sheet, so each literal is written with a leading
# apostrophe to force Excel to store it verbatim (keeps trailing zeros,
# the "%" sign, and negative signs intact instead of being renormalised
# as a number).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'294.10"
$ws.Range("E2").Value = "'1.12%"
$ws.Range("D3").Value = "'31.15"
$ws.Range("E3").Value = "'0.93%"
$ws.Range("D4").Value = "'4.924"
$ws.Range("E4").Value = "'-0.33%"
$ws.Range("D5").Value = "'0.07402"
$ws.Range("E5").Value = "'3.20%"
$ws.Range("D6").Value = "'2.237"
$ws.Range("E6").Value = "'25.14%"
$ws.Range("D7").Value = "'7.754"
$ws.Range("E7").Value = "'1.27%"
$ws.Range("E8").Value = "'0.06%"
$ws.Range("D9").Value = "'0.9153"
$ws.Range("E9").Value = "'2.33%"
$ws.Range("D10").Value = "'0.09146"
$ws.Range("E10").Value = "'19.12%"
$ws.Range("D11").Value = "'0.1700"
$ws.Range("E11").Value = "'2.96%"
$ws.Range("D12").Value = "'0.08320"
$ws.Range("E12").Value = "'3.53%"
$ws.Range("D13").Value = "'0.03126"
$ws.Range("E13").Value = "'2.17%"
$ws.Range("D14").Value = "'0.09974"
$ws.Range("E14").Value = "'-0.53%"
$ws.Range("D15").Value = "'0.001521"
$ws.Range("E15").Value = "'0.90%"
$ws.Range("D16").Value = "'0.005778"
$ws.Range("E16").Value = "'-0.59%"
$ws.Range("E17").Value = "'0.91%"
$ws.Range("D18").Value = "'2.078"
$ws.Range("E18").Value = "'-0.09%"
$ws.Range("E19").Value = "'1.50%"
$ws.Range("D20").Value = "'0.1292"
$ws.Range("E20").Value = "'1.52%"
$ws.Range("D21").Value = "'3.991"
$ws.Range("E21").Value = "'-1.26%"
$ws.Range("E22").Value = "'5.21%"
$ws.Range("D23").Value = "'0.04562"
$ws.Range("E23").Value = "'1.05%"
$ws.Range("E24").Value = "'0.19%"
$ws.Range("D25").Value = "'0.004598"
$ws.Range("E25").Value = "'14.77%"
$ws.Range("D26").Value = "'0.0001305"
$ws.Range("E26").Value = "'4.44%"
$ws.Range("D27").Value = "'0.0003399"
$ws.Range("D39").Value = "'0.01604"
$ws.Range("E39").Value = "'0.32%"
$ws.Range("D40").Value = "'0.04525"
$ws.Range("E40").Value = "'3.30%"
$ws.Range("D41").Value = "'0.007336"
$ws.Range("E41").Value = "'-0.01%"
$ws.Range("D42").Value = "'0.009844"
$ws.Range("E42").Value = "'28.61%"
$ws.Range("D43").Value = "'0.1328"
$ws.Range("E43").Value = "'1.66%"
$ws.Range("D44").Value = "'0.002233"
$ws.Range("E44").Value = "'9.02%"
$ws.Range("D45").Value = "'0.009160"
$ws.Range("E45").Value = "'-0.90%"
$ws.Range("D46").Value = "'0.00006089"
$ws.Range("E46").Value = "'2.32%"
$ws.Range("E47").Value = "'0.18%"
$ws.Range("D48").Value = "'2.282"
$ws.Range("E48").Value = "'1.64%"
$ws.Range("D49").Value = "'0.002003"
$ws.Range("E49").Value = "'-33.23%"
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("E50").Value = "'0.18%"
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("E51").Value = "'0.18%"
